$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.25

# Row 3
$ws.Range("B3").Value = 1.59
$ws.Range("D3").Value = 1.38
$ws.Range("E3").Value = 1.3

# Row 5
$ws.Range("C5").Value = 1.38
$ws.Range("D5").Value = 1.33
$ws.Range("F5").Value = 1.05
$ws.Range("G5").Value = 0.75

# Row 6
$ws.Range("G6").Value = 0.9399999999999999

# Row 7
$ws.Range("E7").Value = 1.88
$ws.Range("F7").Value = 1.5
